$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($r, $c, $v) {
    $ws.Cells.Item($r, $c).Value2 = $v
}

# Remove all hyperlinks (Google Sheets export style: plain-text URLs, no rich hyperlinks)
$ws.Hyperlinks.Delete()

# Remove the now-unused "Hyperlink" cell style from column B (B2:B26)
$ws.Range("B2:B26").Style = "Normal"

Set-Cell 27 1 '2025 Technology Summer Internship - Early Careers (Software or ...'
Set-Cell 27 2 'https://www.wellsfargojobs.com/en/jobs/r-385829/2025-technology-summer-internship-early-careers-software-or-data-engineer/'
Set-Cell 28 1 'SimplifyJobs/Summer2025-Internships: Collection of ... - GitHub'
Set-Cell 28 2 'https://github.com/SimplifyJobs/Summer2025-Internships'
Set-Cell 29 1 'Internships | Museum of Fine Arts Boston'
Set-Cell 29 2 'https://www.mfa.org/working-at-the-mfa/internships'
Set-Cell 30 1 'Data and Analytics jobs in Boston, MA | The Muse'
Set-Cell 30 2 'https://www.themuse.com/hiring/location/boston-ma/category/data_analytics/'
Set-Cell 31 1 'UMass Cybersecurity Institute'
Set-Cell 31 2 'https://infosec.cs.umass.edu/'
Set-Cell 32 1 'Internships - Red Hat Research'
Set-Cell 32 2 'https://research.redhat.com/internships/'
Set-Cell 33 1 'Master of Science in Electrical and Computer Engineering : College ...'
Set-Cell 33 2 'https://www.umass.edu/engineering/academics/ms-electrical-computer-engineering'
Set-Cell 34 1 'Master of Science in Software Development | BU MET'
Set-Cell 34 2 'https://www.bu.edu/met/degrees-certificates/ms-software-development/'
Set-Cell 35 1 'STUDENT AND GRADUATE OPPORTUNITIES'
Set-Cell 35 2 'https://www.santandercareers.com/students'
Set-Cell 36 1 'Internship Programmes'
Set-Cell 36 2 'https://search.jobs.barclays/internships'
Set-Cell 37 1 'Early Career & Internships | Datadog Careers'
Set-Cell 37 2 'https://careers.datadoghq.com/early-careers/'
Set-Cell 38 1 'Internships | Tesla'
Set-Cell 38 2 'https://www.tesla.com/careers/internships'
Set-Cell 39 1 'Wayfair Careers | Students'
Set-Cell 39 2 'https://www.aboutwayfair.com/careers/us-students'
Set-Cell 40 1 'Explore Internships at Procter & Gamble'
Set-Cell 40 2 'https://www.pgcareers.com/global/en/internships'
Set-Cell 41 1 'Students | BNY'
Set-Cell 41 2 'https://www.bnymellon.com/us/en/careers/students.html'
Set-Cell 42 1 'Student Programs & Early Careers | Siemens Software'
Set-Cell 42 2 'https://www.sw.siemens.com/en-US/careers/student-programs-and-early-careers/'
Set-Cell 43 1 'Internships - Students and Graduates | Careers | Oracle'
Set-Cell 43 2 'https://www.oracle.com/careers/students-grads/internships/'
Set-Cell 44 1 'Students | Fidelity Careers'
Set-Cell 44 2 'https://jobs.fidelity.com/students/'
Set-Cell 45 1 'Careers at BlackRock'
Set-Cell 45 2 'https://careers.blackrock.com/'
Set-Cell 46 1 'Internships & Programs | Bain & Company'
Set-Cell 46 2 'https://www.bain.com/careers/work-with-us/internships-programs/'
Set-Cell 47 1 'Campus Reach | Southwest Careers'
Set-Cell 47 2 'https://careers.southwestair.com/campus-reach'
Set-Cell 48 1 'Students | Life at Spotify'
Set-Cell 48 2 'https://www.lifeatspotify.com/students'
Set-Cell 49 1 'Internships - IBM Careers'
Set-Cell 49 2 'https://www.ibm.com/careers/internships'
Set-Cell 50 1 'Early in profession | Microsoft Careers'
Set-Cell 50 2 'https://careers.microsoft.com/v2/global/en/students'
Set-Cell 51 1 'Internships – The Estée Lauder Companies Inc.'
Set-Cell 51 2 'https://www.elcompanies.com/en/careers/students/internships'
Set-Cell 52 1 'Summer Internship Opportunities'
Set-Cell 52 2 'https://jobs.td.com/en/campus-recruitment/summer-internship-opportunities/'

# Column widths (A: 67 -> 72, B: 119 -> 124), offset by Excel's ~0.83 padding quirk
$ws.Columns.Item(1).ColumnWidth = 71.17
$ws.Columns.Item(2).ColumnWidth = 123.17

# Page margins (inches): left/right 0.7 -> 0.75, top/bottom 0.75 -> 1, header/footer 0.3 -> 0.5
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Mark a (harmless, password-less) workbook protection flag so a <workbookProtection/> element is emitted
$wb.Protect($null, $true, $false)
